$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": status text change + widened status columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Sheet "zh-cn": status text change, widened columns, handback info filled in
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsZh.Range("I2").Value = "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md"
$wsZh.Range("J2").Value = "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.db1e2ee57c6cae7909fea87bc09230d401e6fd61.zh-cn.xlf"

$wsZh.Range("I3").Value = "efa5ab67-3930-426d-96da-3800b1dc7f0f.md"
$wsZh.Range("J3").Value = "efa5ab67-3930-426d-96da-3800b1dc7f0f.952a633edff5b4c15815d1edf2015b1d148809b1.zh-cn.xlf"

# Recreate every hyperlink on this sheet so the relationship ids interleave in
# cell order (A2, I2, A3, I3), matching how the report generator lays them out.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md", "", "", "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md", "", "", "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/efa5ab67-3930-426d-96da-3800b1dc7f0f.md", "", "", "efa5ab67-3930-426d-96da-3800b1dc7f0f.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/efa5ab67-3930-426d-96da-3800b1dc7f0f.md", "", "", "efa5ab67-3930-426d-96da-3800b1dc7f0f.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": status text change, widened columns, handback completed
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Range("I2").Value = "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md"
$wsDe.Range("J2").Value = "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.db1e2ee57c6cae7909fea87bc09230d401e6fd61.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-18 13:13:54"

$wsDe.Range("I3").Value = "efa5ab67-3930-426d-96da-3800b1dc7f0f.md"
$wsDe.Range("J3").Value = "efa5ab67-3930-426d-96da-3800b1dc7f0f.952a633edff5b4c15815d1edf2015b1d148809b1.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-18 13:13:54"

# Recreate every hyperlink on this sheet so the relationship ids interleave in
# cell order (A2, I2, A3, I3), matching how the report generator lays them out.
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md", "", "", "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md", "", "", "bf26a85a-a0d0-4a14-b02f-de4b47e5e35b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/efa5ab67-3930-426d-96da-3800b1dc7f0f.md", "", "", "efa5ab67-3930-426d-96da-3800b1dc7f0f.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/703b9969ce072cf492adff56273ae15ae077c1ab/e2e/efa5ab67-3930-426d-96da-3800b1dc7f0f.md", "", "", "efa5ab67-3930-426d-96da-3800b1dc7f0f.md") | Out-Null
